$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values cell by cell, in sheet order.
# A handful of the new D-column prices are plain decimals (e.g. "0.506"); Excel would
# auto-convert those to numeric cells on assignment, but the source keeps every Price/
# Volume cell as text, so we force a text format first and restore the default "Normal"
# style afterwards (keeps the cell text-typed without leaving a stray quote-prefix style).
$ws.Range('D2').Value = '68.408.52'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '2.450.12'
$ws.Range('E3').Value = '  -2.31%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '562.69'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.73'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.37%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.506'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('E9').Value = '  -6.25%  '
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.340'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.78%  '
$ws.Range('E12').Value = '  -2.89%  '
$ws.Range('D13').Value = '2.905.49'
$ws.Range('E13').Value = '  -1.59%  '
$ws.Range('D14').Value = '68.309.05'
$ws.Range('E14').Value = '  -1.92%  '
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '23.60'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.18%  '
$ws.Range('D17').Value = '2.466.88'
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('E18').Value = '  -2.33%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '345.14'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.15'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.57%  '
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('E22').Value = '  -3.57%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.21'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.30%  '
$ws.Range('E25').Value = '  -5.59%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.580.14'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.02'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.26'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -6.35%  '
$ws.Range('E29').Value = '  -6.09%  '
$ws.Range('E30').Value = '  -6.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '437.00'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -5.33%  '
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('B34').Value = 'POPCAT'
$ws.Range('C34').Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.23'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +114.50%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.68'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.19%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '156.36'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('E39').Value = '  -6.18%  '
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('E41').Value = '  -4.02%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.49'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.45%  '
$ws.Range('E43').Value = '  -4.54%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.12'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.72%  '
$ws.Range('E45').Value = '  -5.54%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '135.34'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('E47').Value = '  -3.36%  '
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.486'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -6.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.564'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0915'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.63%  '
